$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shp = $s.Shapes.Item(2)

# Reposition / resize the content placeholder (it had no explicit
# spPr before, now gets an explicit xfrm slightly lower than the
# inherited layout position).
$shp.Left = 43.375040370078736
$shp.Top = 183.39370078740157
$shp.Width = 873.2499212598425
$shp.Height = 313.3563082125984

# Fill in the four paragraphs of "risker/möjligheter" bullet text.
$tr = $shp.TextFrame.TextRange
$tr.Text = "identifierade både projekt- och produktrisker"
$tr.LanguageID = "sv-SE"

$r2 = $tr.InsertAfter([char]13 + "bristfällig kommunikation bland projektmedlemmarna.")
$r2.LanguageID = "sv-SE"

$r3 = $tr.InsertAfter([char]13 + "produkten inte skulle vara tillräckligt intuitiv")
$r3.LanguageID = "sv-SE"

$r4 = $tr.InsertAfter([char]13 + "Möjligheter " + [char]8211 + "> bättre planering ")
$r4.LanguageID = "sv-SE"
